$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'248.54"
$ws.Range("D3").Value = "'22.75"
$ws.Range("D4").Value = "'5.266"
$ws.Range("D5").Value = "'0.05690"
$ws.Range("D6").Value = "'3.419"
$ws.Range("D7").Value = "'6.324"
$ws.Range("D8").Value = "'0.8062"
$ws.Range("D9").Value = "'0.8884"
$ws.Range("B10").Value = "One"
$ws.Range("C10").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D10").Value = "'0.0005813"
$ws.Range("E10").Value = "9OneONE"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "'0.1424"
$ws.Range("E11").Value = "10WazirXWRX"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.07430"
$ws.Range("E12").Value = "11MandalaExchangeTokenMDX"
$ws.Range("B13").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C13").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D13").Value = "'0.03059"
$ws.Range("E13").Value = "12LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D14").Value = "'0.03096"
$ws.Range("E14").Value = "13BitrueCoinBTR"
$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D15").Value = "'0.09394"
$ws.Range("E15").Value = "14BitMartTokenBMX"
$ws.Range("B16").Value = "MCDex"
$ws.Range("C16").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D16").Value = "'3.859"
$ws.Range("E16").Value = "15MCDexMCB"
$ws.Range("B17").Value = "BitForexToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D17").Value = "'0.001589"
$ws.Range("E17").Value = "16BitForexTokenBF"
$ws.Range("B18").Value = "CoinExToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D18").Value = "'0.04796"
$ws.Range("E18").Value = "17CoinExTokenCET"
$ws.Range("B19").Value = "UpBots"
$ws.Range("C19").Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$ws.Range("D19").Value = "'0.01829"
$ws.Range("E19").Value = "18UpBotsUBXTBestin24h"
$ws.Range("D20").Value = "'0.006432"
$ws.Range("D21").Value = "'0.004977"
$ws.Range("D22").Value = "'0.0009968"
$ws.Range("D23").Value = "'0.0001501"
$ws.Range("D24").Value = "'3.688"
$ws.Range("D25").Value = "'2.166"
$ws.Range("D26").Value = "'0.3258"
$ws.Range("D27").Value = "'0.1370"
$ws.Range("D40").Value = "'0.03962"
$ws.Range("D41").Value = "'0.006804"
$ws.Range("D42").Value = "'0.1069"
$ws.Range("D43").Value = "'0.003202"
$ws.Range("D44").Value = "'0.007868"
$ws.Range("D45").Value = "'0.00005588"
$ws.Range("D47").Value = "'0.4993"
$ws.Range("D48").Value = "'0.2010"
$ws.Range("D49").Value = "'0.00002101"
$ws.Range("D50").Value = "'0.01011"
